# "Generate Report for Handback"
#
# The localization report previously listed the three tracked files in
# this fixed order on every sheet:
#   row2 = 3068b5cf-36c3-478a-a79a-02d4decd6479
#   row3 = 86449b45-6e9a-475e-a548-922e5e0d8193
#   row4 = d6745a91-71b4-416f-93c0-c571990ec9a6
#
# A new handback run came in for d6745a91-...: the handback transform
# failed (the produced file name didn't match the expected handoff file
# name), so that file now sorts first, and its status / error detail
# reflect the failure. The other two rows keep their ("Ready for
# handoff") data but shift down to make room.

$wb = $excel.ActiveWorkbook

function Set-RowData {
    param($ws, [int]$row, [string]$mdName, [string]$mdUrl, [string]$status, [string]$date)
    $ws.Range("A$row").Value = $mdName
    $ws.Range("B$row").Value = $status
    $ws.Range("C$row").Value = $status
    $ws.Range("D$row").Value = $date
    $ws.Hyperlinks.Add($ws.Range("A$row"), $mdUrl, $null, $null, $mdName) | Out-Null
}

# ---------------------------------------------------------------------
# Markdown-file hyperlink targets (unchanged across the edit — only the
# row each one is displayed on moves).
# ---------------------------------------------------------------------
$mdUrl_3068 = "https://github.com/OpenLocalizationTest/oltest/blob/e7f6d9994f2051e4e5b6632d66b08c11d6e148ac/e2e/3068b5cf-36c3-478a-a79a-02d4decd6479.md"
$mdUrl_8644 = "https://github.com/OpenLocalizationTest/oltest/blob/e7f6d9994f2051e4e5b6632d66b08c11d6e148ac/e2e/86449b45-6e9a-475e-a548-922e5e0d8193.md"
$mdUrl_d674 = "https://github.com/OpenLocalizationTest/oltest/blob/cb6606bc02ca861046cbb346ee22132c27005229/e2e/d6745a91-71b4-416f-93c0-c571990ec9a6.md"

$name_3068 = "3068b5cf-36c3-478a-a79a-02d4decd6479.md"
$name_8644 = "86449b45-6e9a-475e-a548-922e5e0d8193.md"
$name_d674 = "d6745a91-71b4-416f-93c0-c571990ec9a6.md"

# =======================================================================
# Sheet "Overview"
# =======================================================================
$ov = $wb.Worksheets.Item("Overview")
$ov.Hyperlinks.Delete()

Set-RowData $ov 2 $name_d674 $mdUrl_d674 "Ready for handoff" "2016-49-13 16:49:51"
Set-RowData $ov 3 $name_3068 $mdUrl_3068 "Ready for handoff" "2016-50-13 16:50:38"
Set-RowData $ov 4 $name_8644 $mdUrl_8644 "Ready for handoff" "2016-49-13 16:49:51"

# =======================================================================
# Per-language detail sheets (zh-cn, de-de)
# =======================================================================

function Set-LangSheet {
    param($ws, [string]$lang)

    $ws.Hyperlinks.Delete()

    # xlf hyperlink targets for this language
    if ($lang -eq "zh-cn") {
        $xlfUrl_3068 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5470b9148651a35570e00f6587938f693abf307a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/3068b5cf-36c3-478a-a79a-02d4decd6479.567303eea8a433e1f12ec5976768969d6befd35b.zh-cn.xlf"
        $xlfUrl_8644 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5470b9148651a35570e00f6587938f693abf307a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/86449b45-6e9a-475e-a548-922e5e0d8193.92369443cf2298bdd350fcb0834862425f0818c0.zh-cn.xlf"
        $xlfUrl_d674 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1d1b6232588b822809fac9589e341d26f215300c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d6745a91-71b4-416f-93c0-c571990ec9a6.12766b8e48f7afdf62a9396375eed44b8931d39c.zh-cn.xlf"

        $xlfName_3068 = "3068b5cf-36c3-478a-a79a-02d4decd6479.567303eea8a433e1f12ec5976768969d6befd35b.zh-cn.xlf"
        $xlfName_8644 = "86449b45-6e9a-475e-a548-922e5e0d8193.92369443cf2298bdd350fcb0834862425f0818c0.zh-cn.xlf"
        $xlfName_d674 = "d6745a91-71b4-416f-93c0-c571990ec9a6.12766b8e48f7afdf62a9396375eed44b8931d39c.zh-cn.xlf"

        $handoffDate_3068 = "2016-03-13 16:47:53"
        $handoffDate_8644 = "2016-03-13 16:47:53"
        $handoffDate_d674 = "2016-03-13 16:50:35"
    }
    else {
        $xlfUrl_3068 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ce25f650c1c002d2cafda1d3e373441db12b3a77/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/3068b5cf-36c3-478a-a79a-02d4decd6479.567303eea8a433e1f12ec5976768969d6befd35b.de-de.xlf"
        $xlfUrl_8644 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ce25f650c1c002d2cafda1d3e373441db12b3a77/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/86449b45-6e9a-475e-a548-922e5e0d8193.92369443cf2298bdd350fcb0834862425f0818c0.de-de.xlf"
        $xlfUrl_d674 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/338d97c391374e6242e0d4bc50ca0dc7e164177c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d6745a91-71b4-416f-93c0-c571990ec9a6.12766b8e48f7afdf62a9396375eed44b8931d39c.de-de.xlf"

        $xlfName_3068 = "3068b5cf-36c3-478a-a79a-02d4decd6479.567303eea8a433e1f12ec5976768969d6befd35b.de-de.xlf"
        $xlfName_8644 = "86449b45-6e9a-475e-a548-922e5e0d8193.92369443cf2298bdd350fcb0834862425f0818c0.de-de.xlf"
        $xlfName_d674 = "d6745a91-71b4-416f-93c0-c571990ec9a6.12766b8e48f7afdf62a9396375eed44b8931d39c.de-de.xlf"

        $handoffDate_3068 = "2016-03-13 16:49:51"
        $handoffDate_8644 = "2016-03-13 16:49:51"
        $handoffDate_d674 = "2016-03-13 16:50:38"
    }

    # --- Row 1 : headers (unchanged text, restated so nothing is lost) ---
    $headers = @("Source File Name", "File Extension", "Status", "Latest Handoff File", `
                 "Latest Handoff Datetime", "Latest Target File", "Latest Handback File", `
                 "Latest Handback DateTime", "Handoff Reason", "Dependency From", "Error Detail")
    $col = 1
    foreach ($h in $headers) {
        $ws.Cells.Item(1, $col).Value = $h
        $col++
    }

    # --- Row 2 : d6745a91 — handback transform failed ---
    $ws.Range("A2").Value = $name_d674
    $ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl_d674, $null, $null, $name_d674) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("B2"), $mdUrl_d674, $null, $null, ".md") | Out-Null
    $ws.Range("C2").Value = "Handback transform failed"
    $ws.Hyperlinks.Add($ws.Range("D2"), $xlfUrl_d674, $null, $null, $xlfName_d674) | Out-Null
    $ws.Range("E2").Value = $handoffDate_d674
    $ws.Range("H2").Value = "0001-01-01 00:00:00"
    $ws.Range("I2").Value = "Include"
    $ws.Range("K2").Value = "Handback file name: 5wfi2d3i.4or is different with handoff file name: d6745a91-71b4-416f-93c0-c571990ec9a6.12766b8e48f7afdf62a9396375eed44b8931d39c.$lang."

    # --- Row 3 : 3068b5cf — ready for handoff ---
    $ws.Range("A3").Value = $name_3068
    $ws.Hyperlinks.Add($ws.Range("A3"), $mdUrl_3068, $null, $null, $name_3068) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("B3"), $mdUrl_3068, $null, $null, ".md") | Out-Null
    $ws.Range("C3").Value = "Ready for handoff"
    $ws.Hyperlinks.Add($ws.Range("D3"), $xlfUrl_3068, $null, $null, $xlfName_3068) | Out-Null
    $ws.Range("E3").Value = $handoffDate_3068
    $ws.Range("H3").Value = "0001-01-01 00:00:00"
    $ws.Range("I3").Value = "Include"

    # --- Row 4 : 86449b45 — ready for handoff ---
    $ws.Range("A4").Value = $name_8644
    $ws.Hyperlinks.Add($ws.Range("A4"), $mdUrl_8644, $null, $null, $name_8644) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("B4"), $mdUrl_8644, $null, $null, ".md") | Out-Null
    $ws.Range("C4").Value = "Ready for handoff"
    $ws.Hyperlinks.Add($ws.Range("D4"), $xlfUrl_8644, $null, $null, $xlfName_8644) | Out-Null
    $ws.Range("E4").Value = $handoffDate_8644
    $ws.Range("H4").Value = "0001-01-01 00:00:00"
    $ws.Range("I4").Value = "Include"
}

$zh = $wb.Worksheets.Item("zh-cn")
Set-LangSheet $zh "zh-cn"

$de = $wb.Worksheets.Item("de-de")
Set-LangSheet $de "de-de"

Write-Host "Handback report regenerated."
